$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H6").Value = 77.40000000000001
$ws_ALC.Range("I6").Value = 76.458336
$ws_ALC.Range("J6").Value = 100
$ws_ALC.Range("K6").Value = 229.375008
$ws_ALC.Range("L6").Value = 300
$ws_ALC.Range("M6").Value = -117.375008
$ws_ALC.Range("N6").Value = -524
$ws_ALC.Range("H53").Value = 795.13043
$ws_ALC.Range("I53").Value = 981.46155
$ws_ALC.Range("J53").Value = 552.9
$ws_ALC.Range("K53").Value = 981.46155
$ws_ALC.Range("L53").Value = 552.9
$ws_ALC.Range("M53").Value = -344.46155
$ws_ALC.Range("N53").Value = -1826.9
$ws_ALC.Range("H76").Value = 4666.6665
$ws_ALC.Range("I76").Value = 4500
$ws_ALC.Range("K76").Value = 4500
$ws_ALC.Range("M76").Value = -4185
$ws_ALC.Range("H79").Value = 4666.6665
$ws_ALC.Range("I79").Value = 4500
$ws_ALC.Range("K79").Value = 4500
$ws_ALC.Range("M79").Value = -3408
$ws_ALC.Range("H88").Value = 12890.4
$ws_ALC.Range("I88").Value = 8220.5
$ws_ALC.Range("J88").Value = 14588.546
$ws_ALC.Range("K88").Value = 8220.5
$ws_ALC.Range("L88").Value = 14588.546
$ws_ALC.Range("M88").Value = -7814.5
$ws_ALC.Range("N88").Value = -15400.546
$ws_ALC.Range("H91").Value = 12890.4
$ws_ALC.Range("I91").Value = 8220.5
$ws_ALC.Range("J91").Value = 14588.546
$ws_ALC.Range("K91").Value = 8220.5
$ws_ALC.Range("L91").Value = 14588.546
$ws_ALC.Range("M91").Value = -6816.5
$ws_ALC.Range("N91").Value = -17396.546
$ws_ALC.Range("H99").Value = 2111
$ws_ALC.Range("I99").Value = 1785.7
$ws_ALC.Range("J99").Value = 3737.5
$ws_ALC.Range("K99").Value = 5357.1
$ws_ALC.Range("L99").Value = 11212.5
$ws_ALC.Range("M99").Value = -3859.1
$ws_ALC.Range("N99").Value = -14208.5
$ws_ALC.Range("H101").Value = 2956.1
$ws_ALC.Range("I101").Value = 2760.5
$ws_ALC.Range("K101").Value = 8281.5
$ws_ALC.Range("M101").Value = -6659.5

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 100
$ws_ARM.Range("I5").Value = 100
$ws_ARM.Range("K5").Value = 100
$ws_ARM.Range("M5").Value = 12
$ws_ARM.Range("H74").Value = 3599.3704
$ws_ARM.Range("I74").Value = 3599.3704
$ws_ARM.Range("K74").Value = 3599.3704
$ws_ARM.Range("M74").Value = -2725.3704
$ws_ARM.Range("H77").Value = 3599.3704
$ws_ARM.Range("I77").Value = 3599.3704
$ws_ARM.Range("K77").Value = 17996.852
$ws_ARM.Range("M77").Value = -13628.852
$ws_ARM.Range("H122").Value = 3540.5715
$ws_ARM.Range("I122").Value = 3540.5715
$ws_ARM.Range("K122").Value = 10621.7145
$ws_ARM.Range("M122").Value = -8171.7145

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 100
$ws_BSM.Range("I4").Value = 100
$ws_BSM.Range("K4").Value = 100
$ws_BSM.Range("M4").Value = 15
$ws_BSM.Range("H134").Value = 2433
$ws_BSM.Range("I134").Value = 2319.6
$ws_BSM.Range("J134").Value = 3000
$ws_BSM.Range("K134").Value = 6958.799999999999
$ws_BSM.Range("L134").Value = 9000
$ws_BSM.Range("M134").Value = -4423.799999999999
$ws_BSM.Range("N134").Value = -14070

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1532.1714
$ws_CRP.Range("J31").Value = 1327.6818
$ws_CRP.Range("L31").Value = 1327.6818
$ws_CRP.Range("N31").Value = -1917.6818
$ws_CRP.Range("H34").Value = 1532.1714
$ws_CRP.Range("J34").Value = 1327.6818
$ws_CRP.Range("L34").Value = 1327.6818
$ws_CRP.Range("N34").Value = -1731.6818
$ws_CRP.Range("H62").Value = 9350.299999999999
$ws_CRP.Range("J62").Value = 19201.6
$ws_CRP.Range("L62").Value = 19201.6
$ws_CRP.Range("N62").Value = -20449.6
$ws_CRP.Range("H65").Value = 9350.299999999999
$ws_CRP.Range("J65").Value = 19201.6
$ws_CRP.Range("L65").Value = 96008
$ws_CRP.Range("N65").Value = -102248
$ws_CRP.Range("H107").Value = 2470.75
$ws_CRP.Range("I107").Value = 1838.4
$ws_CRP.Range("K107").Value = 1838.4
$ws_CRP.Range("M107").Value = 81.59999999999991

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 125514.75
$ws_CUL.Range("I4").Value = 21.727272
$ws_CUL.Range("J4").Value = 401599.4
$ws_CUL.Range("K4").Value = 65.181816
$ws_CUL.Range("L4").Value = 1204798.2
$ws_CUL.Range("M4").Value = 46.818184
$ws_CUL.Range("N4").Value = -1205022.2
$ws_CUL.Range("H5").Value = 1518
$ws_CUL.Range("J5").Value = 1448
$ws_CUL.Range("L5").Value = 4344
$ws_CUL.Range("N5").Value = -4568
$ws_CUL.Range("H11").Value = 895.1429000000001
$ws_CUL.Range("I11").Value = 485
$ws_CUL.Range("K11").Value = 1455
$ws_CUL.Range("M11").Value = -1315
$ws_CUL.Range("N118").ClearContents()
$ws_CUL.Range("H118").Value = 1908.75
$ws_CUL.Range("I118").Value = 1908.75
$ws_CUL.Range("J118").Value = 0
$ws_CUL.Range("K118").Value = 5726.25
$ws_CUL.Range("L118").Value = 0
$ws_CUL.Range("M118").Value = -4483.25
$ws_CUL.Range("H135").Value = 1518
$ws_CUL.Range("J135").Value = 1448
$ws_CUL.Range("L135").Value = 13032
$ws_CUL.Range("N135").Value = -18102
$ws_CUL.Range("H140").Value = 2058.476
$ws_CUL.Range("I140").Value = 1543.6842
$ws_CUL.Range("J140").Value = 6949
$ws_CUL.Range("K140").Value = 4631.0526
$ws_CUL.Range("L140").Value = 20847
$ws_CUL.Range("M140").Value = 548.9474
$ws_CUL.Range("N140").Value = -31207

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 176.625
$ws_GSM.Range("I2").Value = 289.875
$ws_GSM.Range("J2").Value = 63.375
$ws_GSM.Range("K2").Value = 289.875
$ws_GSM.Range("L2").Value = 63.375
$ws_GSM.Range("M2").Value = -176.875
$ws_GSM.Range("N2").Value = -289.375
$ws_GSM.Range("H12").Value = 26668
$ws_GSM.Range("J12").Value = 26668
$ws_GSM.Range("L12").Value = 26668
$ws_GSM.Range("N12").Value = -26948
$ws_GSM.Range("N52").ClearContents()
$ws_GSM.Range("H52").Value = 19000
$ws_GSM.Range("J52").Value = 0
$ws_GSM.Range("L52").Value = 0
$ws_GSM.Range("H132").Value = 2709.889
$ws_GSM.Range("I132").Value = 1632.8334
$ws_GSM.Range("K132").Value = 4898.5002
$ws_GSM.Range("M132").Value = -2368.5002

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 2085.4285
$ws_LTW.Range("I7").Value = 2099.6667
$ws_LTW.Range("J7").Value = 2000
$ws_LTW.Range("K7").Value = 2099.6667
$ws_LTW.Range("L7").Value = 2000
$ws_LTW.Range("M7").Value = -1987.6667
$ws_LTW.Range("N7").Value = -2224
$ws_LTW.Range("H74").Value = 22216.875
$ws_LTW.Range("J74").Value = 22216.875
$ws_LTW.Range("L74").Value = 22216.875
$ws_LTW.Range("N74").Value = -24212.875
$ws_LTW.Range("H77").Value = 22216.875
$ws_LTW.Range("J77").Value = 22216.875
$ws_LTW.Range("L77").Value = 66650.625
$ws_LTW.Range("N77").Value = -76634.625
$ws_LTW.Range("H122").Value = 6560.5454
$ws_LTW.Range("I122").Value = 5890.25
$ws_LTW.Range("J122").Value = 6943.5713
$ws_LTW.Range("K122").Value = 17670.75
$ws_LTW.Range("L122").Value = 20830.7139
$ws_LTW.Range("M122").Value = -15220.75
$ws_LTW.Range("N122").Value = -25730.7139
$ws_LTW.Range("H126").Value = 2085.4285
$ws_LTW.Range("I126").Value = 2099.6667
$ws_LTW.Range("J126").Value = 2000
$ws_LTW.Range("K126").Value = 6299.000100000001
$ws_LTW.Range("L126").Value = 6000
$ws_LTW.Range("M126").Value = -3829.000100000001
$ws_LTW.Range("N126").Value = -10940
$ws_LTW.Range("I136").Value = 3039.4
$ws_LTW.Range("K136").Value = 9118.200000000001
$ws_LTW.Range("M136").Value = -6568.200000000001

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H45").Value = 19636.375
$ws_WVR.Range("J45").Value = 23478.4
$ws_WVR.Range("L45").Value = 23478.4
$ws_WVR.Range("N45").Value = -24460.4
$ws_WVR.Range("H113").Value = 2189.6
$ws_WVR.Range("I113").Value = 1987
$ws_WVR.Range("K113").Value = 5961
$ws_WVR.Range("M113").Value = -3791
$ws_WVR.Range("H132").Value = 3032.45
$ws_WVR.Range("I132").Value = 2947
$ws_WVR.Range("K132").Value = 8841
$ws_WVR.Range("M132").Value = -6311
$ws_WVR.Range("H136").Value = 1513.5428
$ws_WVR.Range("I136").Value = 1557.6451
$ws_WVR.Range("J136").Value = 1171.75
$ws_WVR.Range("K136").Value = 4672.9353
$ws_WVR.Range("L136").Value = 3515.25
$ws_WVR.Range("M136").Value = -2122.9353
$ws_WVR.Range("N136").Value = -8615.25

